# Freeze row 1 of anotherExample.xlsx
#
# Also carries the rest of the authored changes that shipped alongside the
# freeze: wider columns, a slightly bigger header font, and two new rows
# at the bottom of the produce table (a merged note row + a "Pear" entry).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths: A=15, B/C/D=20 "characters" ---
# Excel's displayed ColumnWidth and the stored <col width="..."> in the
# XML differ by the default font's padding (~0.8333 chars for Calibri
# 11), so back that out to land on exact stored widths of 15/20/20/20.
$padding = 0.8333333333333333
$ws.Columns.Item(1).ColumnWidth = 15 - $padding
$ws.Columns.Item(2).ColumnWidth = 20 - $padding
$ws.Columns.Item(3).ColumnWidth = 20 - $padding
$ws.Columns.Item(4).ColumnWidth = 20 - $padding

# --- New row 7: a note spanning A7:D7 ---
$ws.Range("A7").Value = "This cell is going to be merged up to D7"
$ws.Range("A7:D7").Merge()

# --- New row 8: another produce entry ---
$ws.Range("A8").Value = "Pear"
$ws.Range("B8").Value = 1.09
$ws.Range("C8").Value = 3.21
$ws.Range("D8").Value = 4.5

# --- Header row formatting: bump the header text to 14pt ---
$ws.Range("A1:D1").Font.Size = 14

# --- Freeze panes so row 1 stays put while scrolling ---
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
[void]$ws.Range("A1").Select()
